# "modified test cases on overdue fix"
#
# - Summary sheet: row 6 (all zeros) is dropped, A5/E5/F5 go from 26.64 to
#   0.97, and an (empty, untouched) G2 cell is introduced so the sheet's
#   used range grows to column G.
# - Repayment schedule sheet: the "Written Off" column (O) is cleared out
#   for rows 3-8 (and the now-orphaned P2 cell is cleared too), and the
#   fee/total figures on rows 3-5 are recomputed to smaller values.
# - Active tab moves from NewLoanInput to Transactions, with each sheet's
#   last remembered selection updated accordingly.

$wb = $excel.ActiveWorkbook

$wsLoan    = $wb.Worksheets.Item("NewLoanInput")
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSched   = $wb.Worksheets.Item("Repayment schedule")
$wsTrans   = $wb.Worksheets.Item("Transactions")

# --- Summary sheet -------------------------------------------------------
# Introduce the (empty) G2 cell that widens the used range to column G.
$wsSummary.Range("G2").Borders.LineStyle = -4142

# Drop the trailing all-zero row.
$wsSummary.Rows.Item(6).Delete()

# Overdue-fee figures shrink from 26.64 to 0.97.
$wsSummary.Range("A5").Value = 0.97
$wsSummary.Range("E5").Value = 0.97
$wsSummary.Range("F5").Value = 0.97

$wsSummary.Range("D5").Select()

# --- Repayment schedule sheet ---------------------------------------------
# The "Written Off" column (O) is no longer populated for rows 3-8, and the
# now-empty P2 cell (row 2 has no "Over Due" total) is removed too.
$wsSched.Range("P2").Clear()
$wsSched.Range("O3:O8").Clear()

$wsSched.Range("J3").Value = 0.51
$wsSched.Range("K3").Value = 888.23
$wsSched.Range("P3").Value = 888.23

$wsSched.Range("J4").Value = 0.46
$wsSched.Range("K4").Value = 888.18
$wsSched.Range("P4").Value = 888.18

$wsSched.Range("J5").Value = 0
$wsSched.Range("K5").Value = 887.72
$wsSched.Range("P5").Value = 887.72

$wsSched.Range("J8").Select()

# --- Transactions sheet ----------------------------------------------------
$wsTrans.Range("G2").Select()

# Transactions becomes the active tab.
$wsTrans.Activate()
